# Fill in the previously-empty "Cucutá" column (O) with its values for each month (rows 3-14)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value  = 22
$ws.Range("O4").Value  = 12
$ws.Range("O5").Value  = 12
$ws.Range("O6").Value  = 43
$ws.Range("O7").Value  = 45
$ws.Range("O8").Value  = 6
$ws.Range("O9").Value  = 7
$ws.Range("O10").Value = 6
$ws.Range("O11").Value = 4
$ws.Range("O12").Value = 2
$ws.Range("O13").Value = 23
$ws.Range("O14").Value = 1

# Update the active selection to reflect where the user ended up (B2) instead of A3:A14
$ws.Range("B2").Select()
